$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H38").Value = 299
$ws_ALC.Range("I38").Value = 254.71428
$ws_ALC.Range("J38").Value = 609
$ws_ALC.Range("K38").Value = 764.14284
$ws_ALC.Range("L38").Value = 1827
$ws_ALC.Range("M38").Value = -392.14284
$ws_ALC.Range("N38").Value = -2571
$ws_ALC.Range("H58").Value = 1310.625
$ws_ALC.Range("J58").Value = 1590.8462
$ws_ALC.Range("L58").Value = 4772.5386
$ws_ALC.Range("N58").Value = -5072.5386
$ws_ALC.Range("H97").Value = 1016.3333
$ws_ALC.Range("J97").Value = 1016.3333
$ws_ALC.Range("L97").Value = 3048.9999
$ws_ALC.Range("N97").Value = -4040.9999
$ws_ALC.Range("H100").Value = 3539.4614
$ws_ALC.Range("I100").Value = 2337
$ws_ALC.Range("J100").Value = 5686.7144
$ws_ALC.Range("K100").Value = 2337
$ws_ALC.Range("L100").Value = 5686.7144
$ws_ALC.Range("M100").Value = -1796
$ws_ALC.Range("N100").Value = -6768.7144
$ws_ALC.Range("H105").Value = 61665.2
$ws_ALC.Range("J105").Value = 61665.2
$ws_ALC.Range("L105").Value = 61665.2
$ws_ALC.Range("N105").Value = -68653.2
$ws_ALC.Range("H113").Value = 9163.115
$ws_ALC.Range("J113").Value = 10547.315
$ws_ALC.Range("L113").Value = 10547.315
$ws_ALC.Range("N113").Value = -17055.315
$ws_ALC.Range("H115").Value = 820.4167
$ws_ALC.Range("I115").Value = 820.4167
$ws_ALC.Range("K115").Value = 2461.2501
$ws_ALC.Range("M115").Value = -894.2501000000002
$ws_ALC.Range("H132").Value = 12719.021
$ws_ALC.Range("I132").Value = 1529.9025
$ws_ALC.Range("J132").Value = 104469.8
$ws_ALC.Range("K132").Value = 4589.7075
$ws_ALC.Range("L132").Value = 313409.4
$ws_ALC.Range("M132").Value = -2059.7075
$ws_ALC.Range("N132").Value = -318469.4
$ws_ALC.Range("H137").Value = 5900.4287
$ws_ALC.Range("I137").Value = 6674.4707
$ws_ALC.Range("K137").Value = 20023.4121
$ws_ALC.Range("M137").Value = -17473.4121

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H63").Value = 6754.8184
$ws_ARM.Range("I63").Value = 4329.143
$ws_ARM.Range("J63").Value = 10999.75
$ws_ARM.Range("K63").Value = 4329.143
$ws_ARM.Range("L63").Value = 10999.75
$ws_ARM.Range("M63").Value = -3643.143
$ws_ARM.Range("N63").Value = -12371.75
$ws_ARM.Range("H66").Value = 6754.8184
$ws_ARM.Range("I66").Value = 4329.143
$ws_ARM.Range("J66").Value = 10999.75
$ws_ARM.Range("K66").Value = 21645.715
$ws_ARM.Range("L66").Value = 54998.75
$ws_ARM.Range("M66").Value = -18213.715
$ws_ARM.Range("N66").Value = -61862.75
$ws_ARM.Range("H97").Value = 2295.7778
$ws_ARM.Range("I97").Value = 1093.25
$ws_ARM.Range("J97").Value = 4700.8335
$ws_ARM.Range("K97").Value = 1093.25
$ws_ARM.Range("L97").Value = 4700.8335
$ws_ARM.Range("M97").Value = -597.25
$ws_ARM.Range("N97").Value = -5692.8335

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H7").Value = 2779.6667
$ws_BSM.Range("I7").Value = 440
$ws_BSM.Range("J7").Value = 3949.5
$ws_BSM.Range("K7").Value = 440
$ws_BSM.Range("L7").Value = 3949.5
$ws_BSM.Range("M7").Value = -327
$ws_BSM.Range("N7").Value = -4175.5
$ws_BSM.Range("H94").Value = 1678.5518
$ws_BSM.Range("I94").Value = 620.5
$ws_BSM.Range("K94").Value = 620.5
$ws_BSM.Range("M94").Value = -169.5

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H5").Value = 964.6667
$ws_CRP.Range("I5").Value = 1001.75
$ws_CRP.Range("J5").Value = 935
$ws_CRP.Range("K5").Value = 1001.75
$ws_CRP.Range("L5").Value = 935
$ws_CRP.Range("M5").Value = -889.75
$ws_CRP.Range("N5").Value = -1159
$ws_CRP.Range("H22").Value = 442.86667
$ws_CRP.Range("I22").Value = 459.6
$ws_CRP.Range("J22").Value = 434.5
$ws_CRP.Range("K22").Value = 459.6
$ws_CRP.Range("L22").Value = 434.5
$ws_CRP.Range("M22").Value = -109.6
$ws_CRP.Range("N22").Value = -1134.5
$ws_CRP.Range("H43").Value = 39088.4
$ws_CRP.Range("J43").Value = 39088.4
$ws_CRP.Range("L43").Value = 39088.4
$ws_CRP.Range("N43").Value = -39456.4
$ws_CRP.Range("H99").Value = 22977730
$ws_CRP.Range("I99").Value = 6912796.5
$ws_CRP.Range("J99").Value = 64287556
$ws_CRP.Range("K99").Value = 6912796.5
$ws_CRP.Range("L99").Value = 64287556
$ws_CRP.Range("M99").Value = -6911298.5
$ws_CRP.Range("N99").Value = -64290552
$ws_CRP.Range("H101").Value = 39088.4
$ws_CRP.Range("J101").Value = 39088.4
$ws_CRP.Range("L101").Value = 39088.4
$ws_CRP.Range("N101").Value = -45578.4
$ws_CRP.Range("H122").Value = 28064702
$ws_CRP.Range("J122").Value = 6299.727
$ws_CRP.Range("L122").Value = 18899.181
$ws_CRP.Range("N122").Value = -23799.181
$ws_CRP.Range("H126").Value = 22977730
$ws_CRP.Range("I126").Value = 6912796.5
$ws_CRP.Range("J126").Value = 64287556
$ws_CRP.Range("K126").Value = 20738389.5
$ws_CRP.Range("L126").Value = 192862668
$ws_CRP.Range("M126").Value = -20735919.5
$ws_CRP.Range("N126").Value = -192867608

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H42").Value = 12900
$ws_CUL.Range("H129").Value = 2247.8125
$ws_CUL.Range("J129").Value = 2849.5
$ws_CUL.Range("L129").Value = 8548.5
$ws_CUL.Range("N129").Value = -18548.5

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H18").Value = 3750
$ws_GSM.Range("I18").Value = 2500
$ws_GSM.Range("K18").Value = 2500
$ws_GSM.Range("M18").Value = -2207
$ws_GSM.Range("H40").Value = 21000
$ws_GSM.Range("I40").Value = 21000
$ws_GSM.Range("J40").Value = 0
$ws_GSM.Range("K40").Value = 21000
$ws_GSM.Range("L40").Value = 0
$ws_GSM.Range("M40").Value = -20849
$ws_GSM.Range("N40").ClearContents()
$ws_GSM.Range("H44").Value = 21833
$ws_GSM.Range("J44").Value = 23666
$ws_GSM.Range("L44").Value = 23666
$ws_GSM.Range("N44").Value = -24858
$ws_GSM.Range("H97").Value = 5170.143
$ws_GSM.Range("J97").Value = 20245.2
$ws_GSM.Range("L97").Value = 20245.2
$ws_GSM.Range("N97").Value = -21237.2
$ws_GSM.Range("H99").Value = 44747.5
$ws_GSM.Range("I99").Value = 0
$ws_GSM.Range("K99").Value = 0
$ws_GSM.Range("M99").ClearContents()
$ws_GSM.Range("H102").Value = 18437046
$ws_GSM.Range("I102").Value = 2859396.5
$ws_GSM.Range("K102").Value = 2859396.5
$ws_GSM.Range("M102").Value = -2857774.5
$ws_GSM.Range("H104").Value = 45000
$ws_GSM.Range("J104").Value = 45000
$ws_GSM.Range("L104").Value = 45000
$ws_GSM.Range("N104").Value = -51988
$ws_GSM.Range("H105").Value = 115000
$ws_GSM.Range("J105").Value = 115000
$ws_GSM.Range("L105").Value = 115000
$ws_GSM.Range("N105").Value = -121988
$ws_GSM.Range("H107").Value = 1063
$ws_GSM.Range("I107").Value = 595
$ws_GSM.Range("J107").Value = 1999
$ws_GSM.Range("K107").Value = 595
$ws_GSM.Range("L107").Value = 1999
$ws_GSM.Range("M107").Value = 1325
$ws_GSM.Range("N107").Value = -5839
$ws_GSM.Range("H122").Value = 3531
$ws_GSM.Range("I122").Value = 3311.25
$ws_GSM.Range("J122").Value = 4410
$ws_GSM.Range("K122").Value = 9933.75
$ws_GSM.Range("L122").Value = 13230
$ws_GSM.Range("M122").Value = -7483.75
$ws_GSM.Range("N122").Value = -18130

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H100").Value = 105245.27
$ws_LTW.Range("I100").Value = 160956.86
$ws_LTW.Range("K100").Value = 160956.86
$ws_LTW.Range("M100").Value = -160415.86

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H46").Value = 73429
$ws_WVR.Range("J46").Value = 73429
$ws_WVR.Range("L46").Value = 73429
$ws_WVR.Range("N46").Value = -73891
$ws_WVR.Range("H81").Value = 4028.2856
$ws_WVR.Range("J81").Value = 4999.8
$ws_WVR.Range("L81").Value = 9999.6
$ws_WVR.Range("N81").Value = -12121.6
$ws_WVR.Range("H84").Value = 4028.2856
$ws_WVR.Range("J84").Value = 4999.8
$ws_WVR.Range("L84").Value = 49998
$ws_WVR.Range("N84").Value = -60606
$ws_WVR.Range("H100").Value = 602.0833
$ws_WVR.Range("I100").Value = 640.375
$ws_WVR.Range("J100").Value = 525.5
$ws_WVR.Range("K100").Value = 1280.75
$ws_WVR.Range("L100").Value = 1051
$ws_WVR.Range("M100").Value = -739.75
$ws_WVR.Range("N100").Value = -2133
$ws_WVR.Range("H122").Value = 1845.6923
$ws_WVR.Range("I122").Value = 1555
$ws_WVR.Range("K122").Value = 4665
$ws_WVR.Range("M122").Value = -2215
$ws_WVR.Range("H134").Value = 73429
$ws_WVR.Range("J134").Value = 73429
$ws_WVR.Range("L134").Value = 220287
$ws_WVR.Range("N134").Value = -225357
